$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: num_customers 72 -> 74 (cohort_size D34 stays 2256); retention_rate recalculated
$ws.Range("C34").Value = 74
$ws.Range("E34").Value = 0.03280141843971631

# Row 36: num_customers 117 -> 118 (cohort_size D36 stays 1930); retention_rate recalculated
$ws.Range("C36").Value = 118
$ws.Range("E36").Value = 0.06113989637305699

# Row 37: num_customers 724 -> 725 and cohort_size 724 -> 725; retention_rate stays 1
$ws.Range("C37").Value = 725
$ws.Range("D37").Value = 725
